# Aggiornate stime con std moving average
$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.00001693170508085783
$wsLinear.Range("B3").Value = -0.0008734629885968446
$wsLinear.Range("B4").Value = 0.00006243000414055855
$wsLinear.Range("B5").Value = "[1.0, 0.3790505108795826, 0.10238248349888494, 0.06095005776014869, 0.044876411349089866, 0.11102549319195869, 0.23375494970670493, 0.2264577542095719, 0.16140099259296792, 0.07905269681729955, 0.03618288716593927, 0.06379614621482893, 0.106193929007722, 0.164932847105379, 0.2525634872153579, 0.16536309393517315, 0.08827020983174376, 0.08035430204502515, 0.031008504069174835, 0.05090244555640475]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.0001993837778804715
$wsNonLinear.Range("B5").Value = -0.0006658084180944305
$wsNonLinear.Range("B6").Value = 0.00005495993789750699
$wsNonLinear.Range("B7").Value = -0.0006187327431428757
$wsNonLinear.Range("B8").Value = 0.0003114010344690967
$wsNonLinear.Range("B9").Value = 0.00007026596836940234
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.3769032131227132, 0.10224003843844251, 0.062203168860321185, 0.04538758058659501, 0.11167710243460031, 0.23244550625748983, 0.2238596332702955, 0.16059620201950528, 0.08006050973051392, 0.036531475387284906, 0.06548908407272853, 0.10740974719202151, 0.1647434492670268, 0.25163405937453365, 0.16530704470576568, 0.09047938040215253, 0.08166493165024083, 0.031701420939380345, 0.05217000868750588]"
